$d = $word.ActiveDocument

# The header table ("Project Name" / "Reviewer's Name" / "Sprint No." /
# "Review Date" / "File Name") is the first table in the document.
# Sprint No. value lives in row 2, column 4; Review Date value lives in
# row 3, column 2 (a merged cell spanning columns 2-4). Editing via the
# Tables/Cell object model (rather than a document-wide Find/Replace)
# keeps the edit scoped to exactly those two cells, so similarly-typed
# text elsewhere in the document (e.g. "1.  License") is left untouched.
$t = $d.Tables.Item(1)

# Sprint No.: "1" -> "2"
$sprintCell = $t.Cell(2, 4)
$sprintRange = $sprintCell.Range
$sprintContent = $d.Range($sprintRange.Start, $sprintRange.End - 1)
$sprintContent.Text = "2"

# Review Date: "02/09/18" -> "02/21/18"
$dateCell = $t.Cell(3, 2)
$dateRange = $dateCell.Range
$dateContent = $d.Range($dateRange.Start, $dateRange.End - 1)
$dateContent.Text = "02/21/18"
